# feat(shifts): process shifts from excel to solver data
#
# Adds a weekly shift-availability table (Monday..Sunday) to the
# "machines" sheet, to the right of the existing Machine ID / Machine
# Name columns. Each machine row gets a default shift window of
# "07:00 - 18:00; 19:00 - 06:00" for Monday through Saturday (Sunday is
# left blank).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("machines")
$ws.Activate()

# --- Header row: day-of-week columns C1:I1 -------------------------------
# Written in this order so new shared-string entries land in the same
# order the source workbook uses.
$ws.Range("C1").Value() = "Monday"
$ws.Range("F1").Value() = "Thursday"
$ws.Range("E1").Value() = "Wednesday"
$ws.Range("D1").Value() = "Tuesday"
$ws.Range("G1").Value() = "Friday"
$ws.Range("H1").Value() = "Saturday"
$ws.Range("I1").Value() = "Sunday"

# --- Data rows: default shift windows for each machine (rows 2-6) -------
$shift = "07:00 - 18:00; 19:00 - 06:00"
$dataRange = $ws.Range("C2:H6")
$dataRange.Value() = $shift
$dataRange.NumberFormat = "h:mm"

# --- Column widths for the new day columns -------------------------------
$ws.Range("C1:I1").ColumnWidth = 23.5

# --- View state: scroll/selection as captured in the source workbook ----
$ws.Range("F11").Select()
$excel.ActiveWindow.ScrollColumn = 2
